$wb = $excel.ActiveWorkbook

# Sheet "展览": F3 1157 -> 1161, F4 2611 -> 2621
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1161
$ws1.Range("F4").Value = 2621

# Sheet "全部类型": F5 1157 -> 1161, F6 2611 -> 2621
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1161
$ws4.Range("F6").Value = 2621
